# Update session numbers in presentations
#
# 1) Bump the "Learning session N" heading on several Course-plan slides.
# 2) Refresh the cached footer date field (09-Apr-22 -> 15-Apr-22) on every
#    slide layout and on the slide master.

$p = $ppt.ActivePresentation

function Set-FirstParagraphText($shape, [string]$newText) {
    # Replace the text of the shape's first paragraph with $newText while
    # keeping it as a single run (selecting the full paragraph's character
    # range rather than assigning TextRange.Text directly avoids leaving
    # behind an extra split run).
    $tr = $shape.TextFrame.TextRange
    $para = $tr.Paragraphs(1, 1)
    $chars = $para.Characters(1, $para.Length)
    $chars.Text = $newText
}

function Set-DatePlaceholderText($shapes, [string]$newText) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $sh = $shapes.Item($k)
        if ($sh.Name -like "Date Placeholder*") {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Length -gt 0) {
                $chars = $tr.Characters(1, $tr.Length)
                $chars.Text = $newText
            }
        }
    }
}

# --- 1) "Learning session N" heading bumps -------------------------------
# slide index -> new heading text
$sessionUpdates = @{
    9  = "Learning session 4"   # was "Learning session 3"
    10 = "Learning session 4"   # was "Learning session 3"
    11 = "Learning session 4"   # was "Learning session 3"
    12 = "Learning session 5"   # was "Learning session 4"
    13 = "Learning session 5"   # was "Learning session 4"
    14 = "Learning session 5"   # was "Learning session 4"
    15 = "Learning session 6"   # was "Learning session 5"
    7  = "Learning session 2-3" # was "Learning session 2"
    8  = "Learning session 2-3" # was "Learning session 2"
}

foreach ($idx in $sessionUpdates.Keys) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item(2)
    Set-FirstParagraphText $shape $sessionUpdates[$idx]
}

# --- 2) Footer date field refresh -----------------------------------------
Set-DatePlaceholderText $p.SlideMaster.Shapes "15-Apr-22"

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Set-DatePlaceholderText $layouts.Item($i).Shapes "15-Apr-22"
}
